$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1, 1, 1, 1, 0, 0, 1, 0, 3, 1, 0, 1, 2, 1, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $values[$i]
}
